$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.840.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.059.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.365"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.829"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +18.67%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.058.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.596"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.16%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.802.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.634.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.074.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000207"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "418.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.236.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.156"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "500.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -14.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0664"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.694"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "154.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.02%  "

# Row reorderings with full B/C/D/E updates (rows 35-39, 43-45)
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.10%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.44%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.131"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.50%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.96%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.356"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.84%  "
